$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 1621.8
$ws.Cells.Item(62, 9).Value = 1621.8
$ws.Cells.Item(62, 11).Value = 1621.8
$ws.Cells.Item(62, 13).Value = -997.8
$ws.Cells.Item(65, 8).Value = 1621.8
$ws.Cells.Item(65, 9).Value = 1621.8
$ws.Cells.Item(65, 11).Value = 8109
$ws.Cells.Item(65, 13).Value = -4989
$ws.Cells.Item(74, 8).Value = 14052.105
$ws.Cells.Item(74, 9).Value = 14528.941
$ws.Cells.Item(74, 11).Value = 14528.941
$ws.Cells.Item(74, 13).Value = -13592.941
$ws.Cells.Item(77, 8).Value = 14052.105
$ws.Cells.Item(77, 9).Value = 14528.941
$ws.Cells.Item(77, 11).Value = 72644.705
$ws.Cells.Item(77, 13).Value = -67964.705
$ws.Cells.Item(118, 8).Value = 1016.44446
$ws.Cells.Item(118, 9).Value = 768.5
$ws.Cells.Item(118, 11).Value = 2305.5
$ws.Cells.Item(118, 13).Value = -648.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2426.8293
$ws.Cells.Item(32, 9).Value = 2426.8293
$ws.Cells.Item(32, 11).Value = 2426.8293
$ws.Cells.Item(32, 13).Value = -2139.8293
$ws.Cells.Item(45, 8).Value = 22265.7
$ws.Cells.Item(45, 9).Value = 31022.785
$ws.Cells.Item(45, 11).Value = 31022.785
$ws.Cells.Item(45, 13).Value = -30645.785
$ws.Cells.Item(122, 8).Value = 2934.389
$ws.Cells.Item(122, 9).Value = 2911.3572
$ws.Cells.Item(122, 11).Value = 8734.071599999999
$ws.Cells.Item(122, 13).Value = -6284.071599999999
$ws.Cells.Item(132, 8).Value = 2660.077
$ws.Cells.Item(132, 9).Value = 2370.2307
$ws.Cells.Item(132, 11).Value = 7110.6921
$ws.Cells.Item(132, 13).Value = -4580.6921

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 25004404
$ws.Cells.Item(20, 9).Value = 29416420
$ws.Cells.Item(20, 11).Value = 29416420
$ws.Cells.Item(20, 13).Value = -29416173
$ws.Cells.Item(134, 8).Value = 2420.366
$ws.Cells.Item(134, 9).Value = 2055.0386
$ws.Cells.Item(134, 10).Value = 3053.6
$ws.Cells.Item(134, 11).Value = 6165.1158
$ws.Cells.Item(134, 12).Value = 9160.799999999999
$ws.Cells.Item(134, 13).Value = -3630.1158
$ws.Cells.Item(134, 14).Value = -14230.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(9, 8).Value = 53890.25
$ws.Cells.Item(9, 10).Value = 53890.25
$ws.Cells.Item(9, 12).Value = 53890.25
$ws.Cells.Item(9, 14).Value = -54226.25
$ws.Cells.Item(58, 8).Value = 2540.476
$ws.Cells.Item(58, 9).Value = 548.3333
$ws.Cells.Item(58, 11).Value = 548.3333
$ws.Cells.Item(58, 13).Value = -345.3333
$ws.Cells.Item(99, 8).Value = 3453.5
$ws.Cells.Item(99, 9).Value = 1999.5
$ws.Cells.Item(99, 10).Value = 4180.5
$ws.Cells.Item(99, 11).Value = 1999.5
$ws.Cells.Item(99, 12).Value = 4180.5
$ws.Cells.Item(99, 13).Value = -501.5
$ws.Cells.Item(99, 14).Value = -7176.5
$ws.Cells.Item(126, 8).Value = 3453.5
$ws.Cells.Item(126, 9).Value = 1999.5
$ws.Cells.Item(126, 10).Value = 4180.5
$ws.Cells.Item(126, 11).Value = 5998.5
$ws.Cells.Item(126, 12).Value = 12541.5
$ws.Cells.Item(126, 13).Value = -3528.5
$ws.Cells.Item(126, 14).Value = -17481.5
$ws.Cells.Item(132, 8).Value = 10105658
$ws.Cells.Item(132, 9).Value = 2073.6667
$ws.Cells.Item(132, 11).Value = 6221.000100000001
$ws.Cells.Item(132, 13).Value = -3691.000100000001
$ws.Cells.Item(136, 8).Value = 2540.476
$ws.Cells.Item(136, 9).Value = 548.3333
$ws.Cells.Item(136, 11).Value = 1644.9999
$ws.Cells.Item(136, 13).Value = 905.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(62, 8).Value = 6000
$ws.Cells.Item(62, 10).Value = 6000
$ws.Cells.Item(62, 12).Value = 18000
$ws.Cells.Item(62, 14).Value = -19372
$ws.Cells.Item(65, 8).Value = 6000
$ws.Cells.Item(65, 10).Value = 6000
$ws.Cells.Item(65, 12).Value = 54000
$ws.Cells.Item(65, 14).Value = -60864
$ws.Cells.Item(68, 8).Value = 9097848
$ws.Cells.Item(68, 9).Value = 799.75
$ws.Cells.Item(68, 11).Value = 2399.25
$ws.Cells.Item(68, 13).Value = -1588.25
$ws.Cells.Item(71, 8).Value = 9097848
$ws.Cells.Item(71, 9).Value = 799.75
$ws.Cells.Item(71, 11).Value = 7197.75
$ws.Cells.Item(71, 13).Value = -3141.75
$ws.Cells.Item(132, 8).Value = 2124.8572
$ws.Cells.Item(132, 9).Value = 1800
$ws.Cells.Item(132, 10).Value = 2179
$ws.Cells.Item(132, 11).Value = 16200
$ws.Cells.Item(132, 12).Value = 19611
$ws.Cells.Item(132, 13).Value = -13670
$ws.Cells.Item(132, 14).Value = -24671

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 6585.8
$ws.Cells.Item(102, 9).Value = 1433.5834
$ws.Cells.Item(102, 11).Value = 1433.5834
$ws.Cells.Item(102, 13).Value = 188.4166
$ws.Cells.Item(126, 8).Value = 8719.286
$ws.Cells.Item(126, 9).Value = 2200.1667
$ws.Cells.Item(126, 11).Value = 6600.500100000001
$ws.Cells.Item(126, 13).Value = -4130.500100000001
$ws.Cells.Item(132, 8).Value = 1932.9678
$ws.Cells.Item(132, 9).Value = 1732.1428
$ws.Cells.Item(132, 11).Value = 5196.428400000001
$ws.Cells.Item(132, 13).Value = -2666.428400000001
$ws.Cells.Item(133, 8).Value = 97418.71000000001
$ws.Cells.Item(133, 9).Value = 96984
$ws.Cells.Item(133, 10).Value = 97491.164
$ws.Cells.Item(133, 11).Value = 96984
$ws.Cells.Item(133, 12).Value = 97491.164
$ws.Cells.Item(133, 13).Value = -91924
$ws.Cells.Item(133, 14).Value = -107611.164
$ws.Cells.Item(136, 8).Value = 116249.75
$ws.Cells.Item(136, 10).Value = 116249.75
$ws.Cells.Item(136, 12).Value = 348749.25
$ws.Cells.Item(136, 14).Value = -353849.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 3126.8462
$ws.Cells.Item(22, 9).Value = 3607
$ws.Cells.Item(22, 10).Value = 2358.6
$ws.Cells.Item(22, 11).Value = 3607
$ws.Cells.Item(22, 12).Value = 2358.6
$ws.Cells.Item(22, 13).Value = -3312
$ws.Cells.Item(22, 14).Value = -2948.6
$ws.Cells.Item(27, 8).Value = 3126.8462
$ws.Cells.Item(27, 9).Value = 3607
$ws.Cells.Item(27, 10).Value = 2358.6
$ws.Cells.Item(27, 11).Value = 3607
$ws.Cells.Item(27, 12).Value = 2358.6
$ws.Cells.Item(27, 13).Value = -3500
$ws.Cells.Item(27, 14).Value = -2572.6
$ws.Cells.Item(68, 8).Value = 3623.4
$ws.Cells.Item(68, 9).Value = 3505.3635
$ws.Cells.Item(68, 10).Value = 3948
$ws.Cells.Item(68, 11).Value = 3505.3635
$ws.Cells.Item(68, 12).Value = 3948
$ws.Cells.Item(68, 13).Value = -2756.3635
$ws.Cells.Item(68, 14).Value = -5446
$ws.Cells.Item(71, 8).Value = 3623.4
$ws.Cells.Item(71, 9).Value = 3505.3635
$ws.Cells.Item(71, 10).Value = 3948
$ws.Cells.Item(71, 11).Value = 17526.8175
$ws.Cells.Item(71, 12).Value = 19740
$ws.Cells.Item(71, 13).Value = -13782.8175
$ws.Cells.Item(71, 14).Value = -27228
$ws.Cells.Item(82, 8).Value = 998.8421
$ws.Cells.Item(82, 10).Value = 920
$ws.Cells.Item(82, 12).Value = 920
$ws.Cells.Item(82, 14).Value = -1642
$ws.Cells.Item(85, 8).Value = 998.8421
$ws.Cells.Item(85, 10).Value = 920
$ws.Cells.Item(85, 12).Value = 920
$ws.Cells.Item(85, 14).Value = -3416
$ws.Cells.Item(93, 8).Value = 1124.8846
$ws.Cells.Item(93, 9).Value = 1125.0454
$ws.Cells.Item(93, 11).Value = 1125.0454
$ws.Cells.Item(93, 13).Value = 122.9546
$ws.Cells.Item(122, 8).Value = 10718.462
$ws.Cells.Item(122, 9).Value = 8032.533
$ws.Cells.Item(122, 10).Value = 14381.091
$ws.Cells.Item(122, 11).Value = 24097.599
$ws.Cells.Item(122, 12).Value = 43143.273
$ws.Cells.Item(122, 13).Value = -21647.599
$ws.Cells.Item(122, 14).Value = -48043.273
$ws.Cells.Item(133, 8).Value = 107775
$ws.Cells.Item(133, 10).Value = 107775
$ws.Cells.Item(133, 12).Value = 107775
$ws.Cells.Item(133, 14).Value = -112835

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 6749.8
$ws.Cells.Item(62, 9).Value = 5583
$ws.Cells.Item(62, 11).Value = 5583
$ws.Cells.Item(62, 13).Value = -4959
$ws.Cells.Item(65, 8).Value = 6749.8
$ws.Cells.Item(65, 9).Value = 5583
$ws.Cells.Item(65, 11).Value = 27915
$ws.Cells.Item(65, 13).Value = -24795
$ws.Cells.Item(116, 8).Value = 0
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 14).Value = -83927.5
